$d = $word.ActiveDocument

# 1. Update the "last generated" timestamp in the Date-styled paragraph.
$d.Content.Find.Execute("June  30, 2021 (05:30:51 PM)", $true, $false, $false, $false, $false,
                         $true, 1, $false, "June  30, 2021 (05:34:59 PM)", 2)

# 2. Promote four top-level section headings from Heading 2 to Heading 1.
#    ("String Comparison" stays Heading 2; "Pushing Further (Optional)" is already Heading 1.)
$headingsToPromote = @(
    "Reading and Understanding",
    "Converting",
    "Comparing",
    "Testing for Equality"
)

foreach ($p in $d.Paragraphs) {
    $text = $p.Range.Text.Trim()
    if ($headingsToPromote -contains $text) {
        $p.Style = "Heading 1"
    }
}
